{"js": "const pairs = [\n  [\"2024-01-19 Friday\", \"2024-01-20 Saturday\"],\n  [\"277\u00f77=39, 4\", \"110\u00f72=55, 0\"],\n  [\"968\u00f75=193, 3\", \"731\u00f74=182, 3\"],\n  [\"966\u00f74=241, 2\", \"219\u00f79=24, 3\"],\n  [\"883\u00f77=126, 1\", \"759\u00f77=108, 3\"],\n  [\"346\u00f79=38, 4\", \"460\u00f77=65, 5\"],\n  [\"381\u00f76=63, 3\", \"461\u00f78=57, 5\"],\n  [\"183\u00f75=36, 3\", \"863\u00f74=215, 3\"],\n  [\"166\u00f78=20, 6\", \"613\u00f78=76, 5\"],\n  [\"298\u00f73=99, 1\", \"575\u00f77=82, 1\"],\n  [\"885\u00f78=110, 5\", \"326\u00f76=54, 2\"],\n  [\"883\u00f79=98, 1\", \"748\u00f77=106, 6\"],\n  [\"522\u00f75=104, 2\", \"962\u00f76=160, 2\"],\n  [\"567\u00f76=94, 3\", \"646\u00f75=129, 1\"],\n  [\"216\u00f76=36, 0\", \"673\u00f79=74, 7\"],\n  [\"951\u00f77=135, 6\", \"293\u00f79=32, 5\"],\n  [\"514\u00f78=64, 2\", \"641\u00f79=71, 2\"],\n  [\"216\u00f77=30, 6\", \"872\u00f74=218, 0\"],\n  [\"360\u00f79=40, 0\", \"154\u00f79=17, 1\"],\n  [\"639\u00f74=159, 3\", \"340\u00f72=170, 0\"],\n  [\"367\u00f75=73, 2\", \"685\u00f73=228, 1\"],\n  [\"277\u00f72=138, 1\", \"946\u00f75=189, 1\"],\n  [\"608\u00f75=121, 3\", \"960\u00f75=192, 0\"],\n  [\"631\u00f74=157, 3\", \"951\u00f73=317, 0\"],\n  [\"133\u00f72=66, 1\", \"701\u00f76=116, 5\"],\n  [\"244\u00f79=27, 1\", \"397\u00f79=44, 1\"],\n];\n\nconst body = context.document.body;\nconst resultSets = pairs.map(([oldText]) => body.search(oldText, { matchCase: true, matchWholeWord: false }));\nawait context.sync();\n\nfor (let i = 0; i < resultSets.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const items = resultSets[i].items;\n  if (items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${items.length}`);\n  }\n  items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2024-01-19 Friday\", \"2024-01-20 Saturday\")\n  ,@(\"277\u00f77=39, 4\", \"110\u00f72=55, 0\")\n  ,@(\"968\u00f75=193, 3\", \"731\u00f74=182, 3\")\n  ,@(\"966\u00f74=241, 2\", \"219\u00f79=24, 3\")\n  ,@(\"883\u00f77=126, 1\", \"759\u00f77=108, 3\")\n  ,@(\"346\u00f79=38, 4\", \"460\u00f77=65, 5\")\n  ,@(\"381\u00f76=63, 3\", \"461\u00f78=57, 5\")\n  ,@(\"183\u00f75=36, 3\", \"863\u00f74=215, 3\")\n  ,@(\"166\u00f78=20, 6\", \"613\u00f78=76, 5\")\n  ,@(\"298\u00f73=99, 1\", \"575\u00f77=82, 1\")\n  ,@(\"885\u00f78=110, 5\", \"326\u00f76=54, 2\")\n  ,@(\"883\u00f79=98, 1\", \"748\u00f77=106, 6\")\n  ,@(\"522\u00f75=104, 2\", \"962\u00f76=160, 2\")\n  ,@(\"567\u00f76=94, 3\", \"646\u00f75=129, 1\")\n  ,@(\"216\u00f76=36, 0\", \"673\u00f79=74, 7\")\n  ,@(\"951\u00f77=135, 6\", \"293\u00f79=32, 5\")\n  ,@(\"514\u00f78=64, 2\", \"641\u00f79=71, 2\")\n  ,@(\"216\u00f77=30, 6\", \"872\u00f74=218, 0\")\n  ,@(\"360\u00f79=40, 0\", \"154\u00f79=17, 1\")\n  ,@(\"639\u00f74=159, 3\", \"340\u00f72=170, 0\")\n  ,@(\"367\u00f75=73, 2\", \"685\u00f73=228, 1\")\n  ,@(\"277\u00f72=138, 1\", \"946\u00f75=189, 1\")\n  ,@(\"608\u00f75=121, 3\", \"960\u00f75=192, 0\")\n  ,@(\"631\u00f74=157, 3\", \"951\u00f73=317, 0\")\n  ,@(\"133\u00f72=66, 1\", \"701\u00f76=116, 5\")\n  ,@(\"244\u00f79=27, 1\", \"397\u00f79=44, 1\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n  if (-not $found) {\n    throw \"Could not find text: $oldText\"\n  }\n}"}
